# Add season record columns (Wins, Losses, Ties) to the sheet.
# This mirrors the commit "Created functions to get season record" which
# adds AD/AE/AF columns containing the team's Wins/Losses/Ties for every
# player row (the same record value for every row, since it's the team's
# season record).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new column headers, matching the bold/bordered
# header style used by the existing header cells (copy formats from A1,
# the same style as every other header cell).
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-49: team's season record (Wins=68, Losses=94, Ties=0)
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 68   # AD
    $ws.Cells.Item($r, 31).Value = 94   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
